$wb = $excel.ActiveWorkbook

# === Sheet: LP1912 ===
$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2, 1).Value = "Última actualización: 16:41:15"
$ws.Cells.Item(3, 1).Value = "Total filas: 137"
$ws.Cells.Item(38, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(39, 3).Value = "14_ABASTO"
$ws.Cells.Item(55, 1).Value = "13:35:25"
$ws.Cells.Item(55, 3).Value = "215C_EL PATO"
$ws.Cells.Item(55, 4).Value = 110
$ws.Cells.Item(57, 1).Value = "13:54:15"
$ws.Cells.Item(57, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(57, 4).Value = 91
$ws.Cells.Item(109, 1).Value = "16:41:15"
$ws.Cells.Item(109, 2).Value = "17:09"
$ws.Cells.Item(109, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(109, 4).Value = 28
$ws.Cells.Item(110, 1).Value = "16:33:53"
$ws.Cells.Item(110, 2).Value = "17:13"
$ws.Cells.Item(110, 4).Value = 40
$ws.Cells.Item(111, 1).Value = "16:22:52"
$ws.Cells.Item(111, 2).Value = "17:14"
$ws.Cells.Item(111, 3).Value = "10_OLMOS"
$ws.Cells.Item(111, 4).Value = 52
$ws.Cells.Item(112, 2).Value = "17:17"
$ws.Cells.Item(112, 3).Value = "17_ROMERO"
$ws.Cells.Item(112, 4).Value = 97
$ws.Cells.Item(113, 1).Value = "15:40:34"
$ws.Cells.Item(113, 2).Value = "17:24"
$ws.Cells.Item(113, 4).Value = 104
$ws.Cells.Item(114, 1).Value = "16:32:08"
$ws.Cells.Item(114, 2).Value = "17:25"
$ws.Cells.Item(114, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(114, 4).Value = 53
$ws.Cells.Item(115, 1).Value = "15:59:48"
$ws.Cells.Item(115, 2).Value = "17:27"
$ws.Cells.Item(115, 3).Value = "15_ABASTO"
$ws.Cells.Item(115, 4).Value = 88
$ws.Cells.Item(116, 1).Value = "16:16:23"
$ws.Cells.Item(116, 2).Value = "17:32"
$ws.Cells.Item(116, 4).Value = 76
$ws.Cells.Item(117, 1).Value = "16:22:52"
$ws.Cells.Item(117, 2).Value = "17:33"
$ws.Cells.Item(117, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(117, 4).Value = 71
$ws.Cells.Item(118, 1).Value = "15:40:34"
$ws.Cells.Item(118, 2).Value = "17:34"
$ws.Cells.Item(118, 3).Value = "10_OLMOS"
$ws.Cells.Item(118, 4).Value = 114
$ws.Cells.Item(119, 1).Value = "16:32:08"
$ws.Cells.Item(119, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(119, 4).Value = 63
$ws.Cells.Item(120, 2).Value = "17:35"
$ws.Cells.Item(120, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(120, 4).Value = 115
$ws.Cells.Item(121, 2).Value = "17:36"
$ws.Cells.Item(121, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(121, 4).Value = 116
$ws.Cells.Item(122, 1).Value = "15:40:34"
$ws.Cells.Item(122, 2).Value = "17:38"
$ws.Cells.Item(122, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(122, 4).Value = 118
$ws.Cells.Item(123, 2).Value = "17:40"
$ws.Cells.Item(123, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(123, 4).Value = 101
$ws.Cells.Item(124, 1).Value = "15:59:48"
$ws.Cells.Item(124, 2).Value = "17:44"
$ws.Cells.Item(124, 4).Value = 105
$ws.Cells.Item(126, 1).Value = "16:32:08"
$ws.Cells.Item(126, 2).Value = "17:45"
$ws.Cells.Item(126, 3).Value = "215B_EL PATO"
$ws.Cells.Item(126, 4).Value = 73
$ws.Cells.Item(127, 1).Value = "16:16:23"
$ws.Cells.Item(127, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(127, 4).Value = 90
$ws.Cells.Item(128, 1).Value = "16:22:52"
$ws.Cells.Item(128, 2).Value = "17:46"
$ws.Cells.Item(128, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(128, 4).Value = 84
$ws.Cells.Item(129, 1).Value = "15:59:48"
$ws.Cells.Item(129, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(129, 4).Value = 108
$ws.Cells.Item(130, 1).Value = "16:33:53"
$ws.Cells.Item(130, 2).Value = "17:47"
$ws.Cells.Item(130, 4).Value = 74
$ws.Cells.Item(131, 1).Value = "16:16:23"
$ws.Cells.Item(131, 2).Value = "17:48"
$ws.Cells.Item(131, 4).Value = 92
$ws.Cells.Item(132, 1).Value = "16:28:39"
$ws.Cells.Item(132, 2).Value = "17:49"
$ws.Cells.Item(132, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(132, 4).Value = 81
$ws.Cells.Item(133, 1).Value = "15:59:48"
$ws.Cells.Item(133, 2).Value = "17:50"
$ws.Cells.Item(133, 4).Value = 111
$ws.Cells.Item(134, 1).Value = "16:16:23"
$ws.Cells.Item(134, 2).Value = "17:51"
$ws.Cells.Item(134, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(134, 4).Value = 95
$ws.Cells.Item(135, 1).Value = "16:22:52"
$ws.Cells.Item(135, 2).Value = "18:02"
$ws.Cells.Item(135, 4).Value = 100
$ws.Cells.Item(136, 2).Value = "18:03"
$ws.Cells.Item(136, 3).Value = "17_ROMERO"
$ws.Cells.Item(136, 4).Value = 107
$ws.Cells.Item(137, 1).Value = "16:16:23"
$ws.Cells.Item(137, 2).Value = "18:04"
$ws.Cells.Item(137, 4).Value = 108
$ws.Cells.Item(138, 1).Value = "16:32:08"
$ws.Cells.Item(138, 2).Value = "18:05"
$ws.Cells.Item(138, 3).Value = "14_ABASTO"
$ws.Cells.Item(138, 4).Value = 93
$ws.Cells.Item(139, 1).Value = "16:28:39"
$ws.Cells.Item(139, 2).Value = "18:24"
$ws.Cells.Item(139, 4).Value = 116
$ws.Cells.Item(140, 1).Value = "16:32:08"
$ws.Cells.Item(140, 2).Value = "18:25"
$ws.Cells.Item(140, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(140, 4).Value = 113
$ws.Cells.Item(140, 5).Value = "LP1912"
$ws.Cells.Item(141, 1).Value = "16:41:15"
$ws.Cells.Item(141, 2).Value = "18:34"
$ws.Cells.Item(141, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(141, 4).Value = 113
$ws.Cells.Item(141, 5).Value = "LP1912"
$ws.Cells.Item(142, 1).Value = "16:41:15"
$ws.Cells.Item(142, 2).Value = "18:38"
$ws.Cells.Item(142, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(142, 4).Value = 117
$ws.Cells.Item(142, 5).Value = "LP1912"

# === Sheet: LP1912-215 ===
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2, 1).Value = "Última actualización: 16:41:15"

# === Sheet: 6203-6173 ===
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2, 1).Value = "Última actualización: 16:41:15"
$ws.Cells.Item(3, 1).Value = "Total filas: 9"
$ws.Cells.Item(8, 1).Value = "16:41:15"
$ws.Cells.Item(8, 2).Value = "17:01"
$ws.Cells.Item(8, 4).Value = 20
$ws.Cells.Item(9, 1).Value = "16:33:53"
$ws.Cells.Item(9, 2).Value = "17:02"
$ws.Cells.Item(9, 4).Value = 29
$ws.Cells.Item(10, 1).Value = "16:32:08"
$ws.Cells.Item(10, 2).Value = "17:03"
$ws.Cells.Item(10, 4).Value = 31
$ws.Cells.Item(11, 1).Value = "16:28:39"
$ws.Cells.Item(11, 2).Value = "17:04"
$ws.Cells.Item(11, 4).Value = 36
$ws.Cells.Item(12, 1).Value = "16:16:23"
$ws.Cells.Item(12, 2).Value = "17:05"
$ws.Cells.Item(12, 4).Value = 49
$ws.Cells.Item(13, 1).Value = "16:22:52"
$ws.Cells.Item(13, 2).Value = "18:21"
$ws.Cells.Item(13, 4).Value = 119
$ws.Cells.Item(14, 1).Value = "16:32:08"
$ws.Cells.Item(14, 2).Value = "18:22"
$ws.Cells.Item(14, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(14, 4).Value = 110
$ws.Cells.Item(14, 5).Value = "L6203"
